# 4.0.3 model and data
# Split the combined "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv"
# input-data file-list rows on the "Boolean" sheet into their per-vehicle-type
# files (LDVs / HDVs / aircraft / rail / ships / motorbikes), and refresh the
# various sheets' navigation state (active sheet/selection) to match the
# author's last-saved view.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Boolean" sheet: expand the two combined trans rows into six rows each
# ---------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" -- make room for five
# more rows so the single row becomes six (rows 17-22), then fill them in.
$wsBool.Rows.Item(17).Resize(5).Insert() | Out-Null

$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# The old "trans/VTQaZ/VTQaZ.csv" row (originally row 21) has shifted down
# by 5 and now sits at row 26. Expand it the same way into six rows.
$wsBool.Rows.Item(26).Resize(5).Insert() | Out-Null

$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of trailing blank formatted rows follow the last data row.
$wsBool.Range("A33:A38").Font.Name = "Calibri"
$wsBool.Range("A33:A38").Font.Size = 11

# ---------------------------------------------------------------------
# Restore each sheet's last-used cursor position / scroll state
# ---------------------------------------------------------------------
$wsBool.Activate() | Out-Null
$wsBool.Range("A32").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10

$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Activate() | Out-Null
$wsInt.Range("A13").Select() | Out-Null

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Range("A1").Select() | Out-Null

Write-Output "done"
